$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.611.70'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.11%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.087.39'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -5.29%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.42'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.22%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.85'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -5.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.073.15'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -5.49%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -8.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.159'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -7.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.43'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.68%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.472'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -6.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.23'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -5.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000228'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -6.44%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.580.46'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.17%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.681.34'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.112'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.087.62'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -5.18%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.77'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -6.90%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '499.05'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -8.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.89'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -8.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.697'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -8.87%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.29'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -7.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.88'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.86%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.46'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -7.19%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.80'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -13.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.01'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.12'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.67%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.08%  '

$ws.Range("B31").Value = 'Stacks'
$ws.Range("C31").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.71'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.57%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.56'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -9.56%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.13'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '526.66'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.92%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.57'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.73%  '

$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.05'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -9.23%  '

$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.55'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0410'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -8.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0812'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.53%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.121'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.66%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.48'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -7.84%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.948.45'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.49%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.69'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.36%  '

$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.19'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.96%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.251'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.41%  '

$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.10%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.64'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.29%  '

$ws.Range("B48").Value = 'PEPE'
$ws.Range("C48").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₃0548'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.03'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -5.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.110'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.10'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -10.69%  '
